# Update the "Förändrad" date column (C) for rows 2-28 from 45182 (2023-09-13)
# to 45184 (2023-09-15), matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
